$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "I believe in open-source, I recently contributed the new reverse proxy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I believe in open-source, I recently contributed to the new reverse proxy",
    2)

$d.Content.Find.Execute(
    "Tech back :", $true, $false, $false, $false, $false, $true, 1, $false,
    "Backend :", 2)

$d.Content.Find.Execute(
    "Tech front :", $true, $false, $false, $false, $false, $true, 1, $false,
    "Frontend :", 2)
